# feat: add 2022-Q4 data
#
# The workbook tracked one quarter of fund-holding detail (sheet "2022-Q1")
# alongside a running "总计" (totals) summary. This adds a 2022-Q4 snapshot:
#   - the existing "2022-Q1" detail sheet is duplicated so its data is
#     preserved unchanged under the same name;
#   - the original detail sheet is overwritten in place with the new
#     2022-Q4 fund-holding rows and renamed "2022-Q4" (ending up
#     positioned - and numbered - ahead of the untouched "2022-Q1" copy);
#   - the "总计" sheet's row 2 is updated to the 2022-Q4 totals, and the
#     former 2022-Q1 totals move down to a new row 3.

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)
$wsQ1 = $wb.Worksheets.Item(2)

# Helper: write $value into $range as literal text, bypassing Excel's
# "looks like a number" auto-conversion, while leaving the cell on the
# default (general) style -- matches cells such as B2/"001534" that must
# stay text-typed without picking up a distinct number-format style.
# A scratch formula (="literal") is computed to a text result, then
# pasted as a value so the destination cell ends up a plain text literal
# with no leftover formula and no new number-format style.
function Set-TextValue {
    param($range, $value)

    $scratch = $wsTotal.Range("ZZ1")
    $escaped = $value -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $wsTotal.Calculate()
    $scratch.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $scratch.ClearContents() | Out-Null
}

# --- Preserve the existing 2022-Q1 detail by duplicating the sheet ---
$wsQ1.Copy($null, $wsQ1) | Out-Null
$wsQ1Copy = $wb.Worksheets.Item($wsQ1.Index + 1)

# --- Turn the original sheet into the new "2022-Q4" detail sheet ---
# (renamed first so the freshly-made copy can reclaim the "2022-Q1" name)
$wsQ4 = $wsQ1
$wsQ4.Name = "2022-Q4"
$wsQ1Copy.Name = "2022-Q1"

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("C2").Value = "华宝万物互联灵活配置混合A"
$wsQ4.Range("H2").Value = 8
Set-TextValue $wsQ4.Range("B2") "001534"
Set-TextValue $wsQ4.Range("D2") "0.77"
Set-TextValue $wsQ4.Range("E2") "91.08"
Set-TextValue $wsQ4.Range("F2") "2.94"
Set-TextValue $wsQ4.Range("G2") "0.0226"

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("C3").Value = "华宝万物互联灵活配置混合C"
$wsQ4.Range("G3").Value = 0
$wsQ4.Range("H3").Value = 8
Set-TextValue $wsQ4.Range("B3") "016463"
Set-TextValue $wsQ4.Range("D3") "0.00"
Set-TextValue $wsQ4.Range("E3") "91.08"
Set-TextValue $wsQ4.Range("F3") "2.94"

# Match the existing "bold header / bordered" look used on the "总计" sheet
# (the 2022-Q4 sheet is new data, so it is restyled from the totals sheet
# rather than inheriting 2022-Q1's header/index-column style).
$wsTotal.Range("B1:D1").Copy() | Out-Null
$wsQ4.Range("B1:H1").PasteSpecial(-4122) | Out-Null
$wsTotal.Range("A2").Copy() | Out-Null
$wsQ4.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# --- Update the "总计" (totals) sheet ---
# Existing row 2 now describes 2022-Q4 instead of 2022-Q1.
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("D2").Value = 0.02

# Re-add the original 2022-Q1 total as a new row 3.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q1"
$wsTotal.Range("C3").Value = 2
$wsTotal.Range("D3").Value = 0.51

$wsTotal.Range("A2").Copy() | Out-Null
$wsTotal.Range("A3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
